$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly to fit the new, longer "Knob" section header text
$ws.Columns.Item(1).ColumnWidth = 19.42578125

# Add a new row (18) below the existing "Silicone Caulk" row (17), mirroring
# the style used by the other "Other:" sub-items (A13, A15, A17 - bold font).
$ws.Range("A18").Value = "Knob"
$ws.Range("A18").Style = $ws.Range("A17").Style

$ws.Range("B18").Value = "1ea"
$ws.Range("C18").Value = "6mm shaft diameter"
$ws.Range("D18").Value = "Source from Ebay, Amazon, Aliexpress. "

# Move selection to reflect where the user left off after adding the row
$ws.Range("A19").Select()
